# "added 4wk low sales check"
# Update the Forecast Comparison sheet (MyForecast, Inventory Coverage,
# Stockout Risk, Reorder Urgency, Seasonality Index) and the dependent
# totals on the Summary sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row => MyForecast (D), Inventory Coverage (H), Stockout Risk (I),
#         Reorder Urgency (J), Seasonality Index (L)
$rows = @(
    @{ Row = 2;  D = 197; H = 10.36;              I = $null;  J = $null;    L = 0.99 },
    @{ Row = 3;  D = 199; H = 9.279999999999999;  I = $null;  J = $null;    L = 0.89 },
    @{ Row = 4;  D = 200; H = 8.199999999999999;  I = $null;  J = $null;    L = 0.8100000000000001 },
    @{ Row = 5;  D = 202; H = 7.14;               I = $null;  J = $null;    L = 1.18 },
    @{ Row = 6;  D = 204; H = 6.09;               I = $null;  J = $null;    L = 1.12 },
    @{ Row = 7;  D = 206; H = 5.05;               I = $null;  J = $null;    L = 1.15 },
    @{ Row = 8;  D = 207; H = 4.02;               I = $null;  J = $null;    L = 0.93 },
    @{ Row = 9;  D = 209; H = 2.99;               I = "Low";  J = "Normal"; L = 1.16 },
    @{ Row = 10; D = 211; H = 1.98;               I = "Low";  J = "Normal"; L = 0.92 },
    @{ Row = 11; D = 212; H = 0.97;               I = "Low";  J = $null;    L = 1 },
    @{ Row = 12; D = 214; H = $null;              I = $null;  J = $null;    L = 1.09 },
    @{ Row = 13; D = 216; H = $null;              I = $null;  J = $null;    L = 1.07 },
    @{ Row = 14; D = 217; H = $null;              I = $null;  J = $null;    L = 1.19 },
    @{ Row = 15; D = 219; H = $null;              I = $null;  J = $null;    L = 0.98 },
    @{ Row = 16; D = 221; H = $null;              I = $null;  J = $null;    L = 1.16 },
    @{ Row = 17; D = 223; H = $null;              I = $null;  J = $null;    L = 0.84 }
)

foreach ($r in $rows) {
    $wsForecast.Range("D" + $r.Row).Value = $r.D
    if ($null -ne $r.H) { $wsForecast.Range("H" + $r.Row).Value = $r.H }
    if ($null -ne $r.I) { $wsForecast.Range("I" + $r.Row).Value = $r.I }
    if ($null -ne $r.J) { $wsForecast.Range("J" + $r.Row).Value = $r.J }
    $wsForecast.Range("L" + $r.Row).Value = $r.L
}

# Summary sheet totals recalculated from the refreshed forecast data.
# Leading apostrophe forces these numeric-looking values to stay text,
# matching the existing (inline string) cell type on this sheet.
$wsSummary.Range("B9").Value  = "'3364"
$wsSummary.Range("B10").Value = "'1627"
$wsSummary.Range("B11").Value = "'800"
$wsSummary.Range("B12").Value = "'224"
$wsSummary.Range("B14").Value = "'197"
